$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MergeData")

# Copy the formatting (style) of CR2 onto the rest of the CR column (CR3:CR57) so every
# cell shares the same border style - CR33 previously carried a different left-over style
# from an older table boundary, and CR3:CR32/CR34:CR57 had no cell/style at all yet.
$ws.Cells.Item(2, 96).Copy($ws.Range("CR3:CR57"))

# CR2 gets its own (non-shared) formula.
$ws.Cells.Item(2, 96).Formula = '="20.12.2025 22:00"'

# CR3:CR57 share one formula - the "son teslim tarihi" (final submission date/time) stamp.
$ws.Range("CR3:CR57").Formula = '="20.12.2025 22:00"'

# Widen column CR (96) so the new date/time text fits.
$ws.Columns.Item(96).ColumnWidth = 14.1666666666667
